$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 17.61317913292495
$ws.Cells.Item(2, 4).Value = 7.888452948504171
$ws.Cells.Item(2, 5).Value = 14.3414213402736
$ws.Cells.Item(2, 6).Value = 42.13282251075115
$ws.Cells.Item(2, 7).Value = 49.80736577915413
$ws.Cells.Item(2, 8).Value = 19.55169442134907
$ws.Cells.Item(2, 10).Value = 11.31980718705681
$ws.Cells.Item(2, 11).Value = 12.02629931463118
$ws.Cells.Item(2, 12).Value = 9.808381591485189
$ws.Cells.Item(2, 13).Value = 16.56971442282067
$ws.Cells.Item(2, 14).Value = 23.30224139801533
# Row 3
$ws.Cells.Item(3, 2).Value = 17.55005965566573
$ws.Cells.Item(3, 4).Value = 7.893639737159547
$ws.Cells.Item(3, 5).Value = 14.35691861751172
$ws.Cells.Item(3, 6).Value = 42.18121235685641
$ws.Cells.Item(3, 7).Value = 49.80588556180674
$ws.Cells.Item(3, 8).Value = 19.59138354901564
$ws.Cells.Item(3, 10).Value = 11.32951137998026
$ws.Cells.Item(3, 11).Value = 11.84940176232971
$ws.Cells.Item(3, 12).Value = 9.797634431228012
$ws.Cells.Item(3, 13).Value = 16.56144917284506
$ws.Cells.Item(3, 14).Value = 23.3660527026018
# Row 4
$ws.Cells.Item(4, 2).Value = 17.51455533801389
$ws.Cells.Item(4, 4).Value = 7.897298070513337
$ws.Cells.Item(4, 5).Value = 14.36704726818232
$ws.Cells.Item(4, 6).Value = 42.21908164985124
$ws.Cells.Item(4, 7).Value = 49.8168129729217
$ws.Cells.Item(4, 8).Value = 19.61876907178497
$ws.Cells.Item(4, 10).Value = 11.33580684458846
$ws.Cells.Item(4, 11).Value = 11.74220416808797
$ws.Cells.Item(4, 12).Value = 9.792449571681122
$ws.Cells.Item(4, 13).Value = 16.55878611947233
$ws.Cells.Item(4, 14).Value = 23.40713058633856
# Row 5
$ws.Cells.Item(5, 2).Value = 17.50091589196503
$ws.Cells.Item(5, 4).Value = 7.898908292424982
$ws.Cells.Item(5, 5).Value = 14.37132933757204
$ws.Cells.Item(5, 6).Value = 42.23656327880951
$ws.Cells.Item(5, 7).Value = 49.82424150727768
$ws.Cells.Item(5, 8).Value = 19.63068692550337
$ws.Cells.Item(5, 10).Value = 11.33845727022461
$ws.Cells.Item(5, 11).Value = 11.69893284302558
$ws.Cells.Item(5, 12).Value = 9.790694219838734
$ws.Cells.Item(5, 13).Value = 16.55830981931069
$ws.Cells.Item(5, 14).Value = 23.42434853614921
# Row 6
$ws.Cells.Item(6, 2).Value = 17.49870144503516
$ws.Cells.Item(6, 4).Value = 7.899182890780988
$ws.Cells.Item(6, 5).Value = 14.372049717512
$ws.Cells.Item(6, 6).Value = 42.23958981741687
$ws.Cells.Item(6, 7).Value = 49.82565466640691
$ws.Cells.Item(6, 8).Value = 19.63271164552031
$ws.Cells.Item(6, 10).Value = 11.33890250996305
$ws.Cells.Item(6, 11).Value = 11.69177421562802
$ws.Cells.Item(6, 12).Value = 9.790424392425491
$ws.Cells.Item(6, 13).Value = 16.5582675733467
$ws.Cells.Item(6, 14).Value = 23.42723649304489
# Row 7
$ws.Cells.Item(7, 2).Value = 17.51436802145205
$ws.Cells.Item(7, 4).Value = 7.897319302555178
$ws.Cells.Item(7, 5).Value = 14.3671043913697
$ws.Cells.Item(7, 6).Value = 42.21930911730522
$ws.Cells.Item(7, 7).Value = 49.81690111137104
$ws.Cells.Item(7, 8).Value = 19.61892673152479
$ws.Cells.Item(7, 10).Value = 11.33584224481038
$ws.Cells.Item(7, 11).Value = 11.7416188518575
$ws.Cells.Item(7, 12).Value = 9.792424448338229
$ws.Cells.Item(7, 13).Value = 16.55877722741642
$ws.Cells.Item(7, 14).Value = 23.40736085506116
# Row 8
$ws.Cells.Item(8, 2).Value = 17.59074776872002
$ws.Cells.Item(8, 4).Value = 7.89014324744721
$ws.Cells.Item(8, 5).Value = 14.34663773149538
$ws.Cells.Item(8, 6).Value = 42.14781330852097
$ws.Cells.Item(8, 7).Value = 49.80439880836029
$ws.Cells.Item(8, 8).Value = 19.56475297317347
$ws.Cells.Item(8, 10).Value = 11.32308336954419
$ws.Cells.Item(8, 11).Value = 11.96504436429949
$ws.Cells.Item(8, 12).Value = 9.804383743362912
$ws.Cells.Item(8, 13).Value = 16.56636551966382
$ws.Cells.Item(8, 14).Value = 23.32385046282541
# Row 9
$ws.Cells.Item(9, 2).Value = 17.76579069331936
$ws.Cells.Item(9, 4).Value = 7.879814985887719
$ws.Cells.Item(9, 5).Value = 14.31135275013036
$ws.Cells.Item(9, 6).Value = 42.07239608969299
$ws.Cells.Item(9, 7).Value = 49.87376279793516
$ws.Cells.Item(9, 8).Value = 19.48246864640856
$ws.Cells.Item(9, 10).Value = 11.30072764281834
$ws.Cells.Item(9, 11).Value = 12.41196591067265
$ws.Cells.Item(9, 12).Value = 9.838961509682084
$ws.Cells.Item(9, 13).Value = 16.60026723638254
$ws.Cells.Item(9, 14).Value = 23.17508493935966
# Row 10
$ws.Cells.Item(10, 2).Value = 17.9090240869461
$ws.Cells.Item(10, 4).Value = 7.874490705179246
$ws.Cells.Item(10, 5).Value = 14.28836374406082
$ws.Cells.Item(10, 6).Value = 42.05653151164692
$ws.Cells.Item(10, 7).Value = 49.98180701473182
$ws.Cells.Item(10, 8).Value = 19.43663727957968
$ws.Cells.Item(10, 10).Value = 11.28591325097134
$ws.Cells.Item(10, 11).Value = 12.74238829804889
$ws.Cells.Item(10, 12).Value = 9.87101343564529
$ws.Cells.Item(10, 13).Value = 16.63660251677673
$ws.Cells.Item(10, 14).Value = 23.07484817284879
# Row 11
$ws.Cells.Item(11, 2).Value = 17.97718055674734
$ws.Cells.Item(11, 4).Value = 7.872556173493184
$ws.Cells.Item(11, 5).Value = 14.27853812670858
$ws.Cells.Item(11, 6).Value = 42.05789859256733
$ws.Cells.Item(11, 7).Value = 50.04328366980329
$ws.Cells.Item(11, 8).Value = 19.41896606155677
$ws.Cells.Item(11, 10).Value = 11.27952052496924
$ws.Cells.Item(11, 11).Value = 12.8924747785135
$ws.Cells.Item(11, 12).Value = 9.887005688504829
$ws.Cells.Item(11, 13).Value = 16.65557274903477
$ws.Cells.Item(11, 14).Value = 23.03119800688608
# Row 12
$ws.Cells.Item(12, 2).Value = 18.0034044929098
$ws.Cells.Item(12, 4).Value = 7.871893379372495
$ws.Cells.Item(12, 5).Value = 14.27490798136948
$ws.Cells.Item(12, 6).Value = 42.05964913274025
$ws.Cells.Item(12, 7).Value = 50.06832658156887
$ws.Cells.Item(12, 8).Value = 19.41273151809983
$ws.Cells.Item(12, 10).Value = 11.27714935868319
$ws.Cells.Item(12, 11).Value = 12.94922143618417
$ws.Cells.Item(12, 12).Value = 9.89326146032322
$ws.Cells.Item(12, 13).Value = 16.66310336492106
$ws.Cells.Item(12, 14).Value = 23.01494767930144
# Row 13
$ws.Cells.Item(13, 2).Value = 17.99773853363777
$ws.Cells.Item(13, 4).Value = 7.872033026190747
$ws.Cells.Item(13, 5).Value = 14.27568577338224
$ws.Cells.Item(13, 6).Value = 42.05921732251283
$ws.Cells.Item(13, 7).Value = 50.0628549033521
$ws.Cells.Item(13, 8).Value = 19.41405390362461
$ws.Cells.Item(13, 10).Value = 11.27765782795645
$ws.Cells.Item(13, 11).Value = 12.93700479982825
$ws.Cells.Item(13, 12).Value = 9.891905331389706
$ws.Cells.Item(13, 13).Value = 16.66146614605831
$ws.Cells.Item(13, 14).Value = 23.01843508468617
# Row 14
$ws.Cells.Item(14, 2).Value = 17.97932980513143
$ws.Cells.Item(14, 4).Value = 7.872500248620715
$ws.Cells.Item(14, 5).Value = 14.27823765839459
$ws.Cells.Item(14, 6).Value = 42.05801791343821
$ws.Cells.Item(14, 7).Value = 50.04530869285851
$ws.Cells.Item(14, 8).Value = 19.41844397876171
$ws.Cells.Item(14, 10).Value = 11.27932445432959
$ws.Cells.Item(14, 11).Value = 12.89714543232188
$ws.Cells.Item(14, 12).Value = 9.887516367914948
$ws.Cells.Item(14, 13).Value = 16.6561853669418
$ws.Cells.Item(14, 14).Value = 23.02985549637453
# Row 15
$ws.Cells.Item(15, 2).Value = 17.96810740953887
$ws.Cells.Item(15, 4).Value = 7.872795512181312
$ws.Cells.Item(15, 5).Value = 14.27981255176877
$ws.Cells.Item(15, 6).Value = 42.05744373820649
$ws.Cells.Item(15, 7).Value = 50.03479041558126
$ws.Cells.Item(15, 8).Value = 19.4211925675844
$ws.Cells.Item(15, 10).Value = 11.28035176741983
$ws.Cells.Item(15, 11).Value = 12.87271736191798
$ws.Cells.Item(15, 12).Value = 9.884853931417492
$ws.Cells.Item(15, 13).Value = 16.65299580289294
$ws.Cells.Item(15, 14).Value = 23.03688713601117
# Row 16
$ws.Cells.Item(16, 2).Value = 17.90462883590573
$ws.Cells.Item(16, 4).Value = 7.874626912995578
$ws.Cells.Item(16, 5).Value = 14.28901856639758
$ws.Cells.Item(16, 6).Value = 42.05661478827719
$ws.Cells.Item(16, 7).Value = 49.97803664590174
$ws.Cells.Item(16, 8).Value = 19.43785609891482
$ws.Cells.Item(16, 10).Value = 11.28633798479758
$ws.Cells.Item(16, 11).Value = 12.73257048765229
$ws.Cells.Item(16, 12).Value = 9.869996446113895
$ws.Cells.Item(16, 13).Value = 16.63541155067314
$ws.Cells.Item(16, 14).Value = 23.0777399205212
# Row 17
$ws.Cells.Item(17, 2).Value = 17.86644321684165
$ws.Cells.Item(17, 4).Value = 7.875875045731124
$ws.Cells.Item(17, 5).Value = 14.2948278596931
$ws.Cells.Item(17, 6).Value = 42.05830396473333
$ws.Cells.Item(17, 7).Value = 49.94637106043407
$ws.Cells.Item(17, 8).Value = 19.44889270008486
$ws.Cells.Item(17, 10).Value = 11.29009892196479
$ws.Cells.Item(17, 11).Value = 12.64649736129978
$ws.Cells.Item(17, 12).Value = 9.861241243548537
$ws.Cells.Item(17, 13).Value = 16.62524662889775
$ws.Cells.Item(17, 14).Value = 23.10329994574393
# Row 18
$ws.Cells.Item(18, 2).Value = 17.84476316881319
$ws.Cells.Item(18, 4).Value = 7.87663884557089
$ws.Cells.Item(18, 5).Value = 14.2982287350179
$ws.Cells.Item(18, 6).Value = 42.06008363170951
$ws.Cells.Item(18, 7).Value = 49.92931880380318
$ws.Cells.Item(18, 8).Value = 19.45553975383469
$ws.Cells.Item(18, 10).Value = 11.29229473124697
$ws.Cells.Item(18, 11).Value = 12.59697248185969
$ws.Cells.Item(18, 12).Value = 9.856338555452158
$ws.Cells.Item(18, 13).Value = 16.61963006344693
$ws.Cells.Item(18, 14).Value = 23.11818481893176
# Row 19
$ws.Cells.Item(19, 2).Value = 17.83747183322062
$ws.Cells.Item(19, 4).Value = 7.876905349687545
$ws.Cells.Item(19, 5).Value = 14.29939044557057
$ws.Cells.Item(19, 6).Value = 42.06082502221772
$ws.Cells.Item(19, 7).Value = 49.9237448562352
$ws.Cells.Item(19, 8).Value = 19.45784169360051
$ws.Cells.Item(19, 10).Value = 11.29304380259739
$ws.Cells.Item(19, 11).Value = 12.58020285802796
$ws.Cells.Item(19, 12).Value = 9.854701537729051
$ws.Cells.Item(19, 13).Value = 16.61776801708181
$ws.Cells.Item(19, 14).Value = 23.1232561175212
# Row 20
$ws.Cells.Item(20, 2).Value = 17.87047893837741
$ws.Cells.Item(20, 4).Value = 7.875737431087182
$ws.Cells.Item(20, 5).Value = 14.29420329229132
$ws.Cells.Item(20, 6).Value = 42.05804052259175
$ws.Cells.Item(20, 7).Value = 49.94962182496195
$ws.Cells.Item(20, 8).Value = 19.44768687823852
$ws.Cells.Item(20, 10).Value = 11.28969518938657
$ws.Cells.Item(20, 11).Value = 12.65566225280018
$ws.Cells.Item(20, 12).Value = 9.862159499465633
$ws.Cells.Item(20, 13).Value = 16.62630492343142
$ws.Cells.Item(20, 14).Value = 23.10056006173285
# Row 21
$ws.Cells.Item(21, 2).Value = 17.98472577993279
$ws.Cells.Item(21, 4).Value = 7.872361123167816
$ws.Cells.Item(21, 5).Value = 14.27748565140006
$ws.Cells.Item(21, 6).Value = 42.058336765115
$ws.Cells.Item(21, 7).Value = 50.05041467750219
$ws.Cells.Item(21, 8).Value = 19.41714209871517
$ws.Cells.Item(21, 10).Value = 11.27883358061709
$ws.Cells.Item(21, 11).Value = 12.9088559123592
$ws.Cells.Item(21, 12).Value = 9.888800114874229
$ws.Cells.Item(21, 13).Value = 16.65772707569537
$ws.Cells.Item(21, 14).Value = 23.0264934800428
# Row 22
$ws.Cells.Item(22, 2).Value = 18.0618008991227
$ws.Cells.Item(22, 4).Value = 7.870561033239795
$ws.Cells.Item(22, 5).Value = 14.26708766912533
$ws.Cells.Item(22, 6).Value = 42.06571520509104
$ws.Cells.Item(22, 7).Value = 50.12655917646369
$ws.Cells.Item(22, 8).Value = 19.39984405540093
$ws.Cells.Item(22, 10).Value = 11.27202402809553
$ws.Cells.Item(22, 11).Value = 13.07379621002819
$ws.Cells.Item(22, 12).Value = 9.907374608208794
$ws.Cells.Item(22, 13).Value = 16.68028380119775
$ws.Cells.Item(22, 14).Value = 22.97971267045635
# Row 23
$ws.Cells.Item(23, 2).Value = 18.0204497769628
$ws.Cells.Item(23, 4).Value = 7.871484688923594
$ws.Cells.Item(23, 5).Value = 14.27258905957928
$ws.Cells.Item(23, 6).Value = 42.06112046206643
$ws.Cells.Item(23, 7).Value = 50.08498331887849
$ws.Cells.Item(23, 8).Value = 19.40883247477021
$ws.Cells.Item(23, 10).Value = 11.27563202188374
$ws.Cells.Item(23, 11).Value = 12.98583142174137
$ws.Cells.Item(23, 12).Value = 9.897355693126306
$ws.Cells.Item(23, 13).Value = 16.66806137220855
$ws.Cells.Item(23, 14).Value = 23.00453204560383
# Row 24
$ws.Cells.Item(24, 2).Value = 17.86865353587339
$ws.Cells.Item(24, 4).Value = 7.875799502657562
$ws.Cells.Item(24, 5).Value = 14.29448546909158
$ws.Cells.Item(24, 6).Value = 42.05815710608486
$ws.Cells.Item(24, 7).Value = 49.94814856293124
$ws.Cells.Item(24, 8).Value = 19.44823108967408
$ws.Cells.Item(24, 10).Value = 11.28987761222578
$ws.Cells.Item(24, 11).Value = 12.6515189265832
$ws.Cells.Item(24, 12).Value = 9.861743948312112
$ws.Cells.Item(24, 13).Value = 16.62582575988454
$ws.Cells.Item(24, 14).Value = 23.10179817110985
# Row 25
$ws.Cells.Item(25, 2).Value = 17.71580879578787
$ws.Cells.Item(25, 4).Value = 7.882210139573881
$ws.Cells.Item(25, 5).Value = 14.32038131299481
$ws.Cells.Item(25, 6).Value = 42.08585500322188
$ws.Cells.Item(25, 7).Value = 49.84495542170176
$ws.Cells.Item(25, 8).Value = 19.50216211210891
$ws.Cells.Item(25, 10).Value = 11.30649166483554
$ws.Cells.Item(25, 11).Value = 12.29046327347665
$ws.Cells.Item(25, 12).Value = 9.828429403931228
$ws.Cells.Item(25, 13).Value = 16.5890758859755
$ws.Cells.Item(25, 14).Value = 23.21373283052421
